$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.797000000000001
$ws.Range("A9").Value = -20.912
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
